$wb = $excel.ActiveWorkbook

# 1. Update the shared "status" string wherever it appears (Overview!B/C, zh-cn!C, de-de!C)
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("B2").Value = "Handed back: in sync with en-US"
$overview.Range("C2").Value = "Handed back: in sync with en-US"
$overview.Range("B3").Value = "Handed back: in sync with en-US"
$overview.Range("C3").Value = "Handed back: in sync with en-US"

$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("C2").Value = "Handed back: in sync with en-US"
$zhcn.Range("C3").Value = "Handed back: in sync with en-US"

$dede = $wb.Worksheets.Item("de-de")
$dede.Range("C2").Value = "Handed back: in sync with en-US"
$dede.Range("C3").Value = "Handed back: in sync with en-US"

# 2. zh-cn: update handback datetime (H2/H3) and add F/G columns (Target File / Handback File)
$zhcn.Range("H2").Value = "2016-03-21 18:46:52"
$zhcn.Range("H3").Value = "2016-03-21 18:46:52"

$zhcn.Range("F2").Value = "748403bf-c93f-4a77-8128-a125cc87e042.md"
$zhcn.Range("G2").Value = "748403bf-c93f-4a77-8128-a125cc87e042.3377a6dd9e3bc8a1d85ec6d6e94610766108e47d.zh-cn.xlf"
$zhcn.Range("F3").Value = "89ec4140-7020-4012-9fe1-624c2b8a2ebb.md"
$zhcn.Range("G3").Value = "89ec4140-7020-4012-9fe1-624c2b8a2ebb.494a8a5e1a4fe5743b433ae8722919d841a9c757.zh-cn.xlf"

$zhcn.Hyperlinks.Add($zhcn.Range("F2"), "https://github.com/OpenLocalizationTest/oltest/blob/0b105620a0bcd802b94d27efa52426285c9e8580/e2e/748403bf-c93f-4a77-8128-a125cc87e042.md", "", "", "748403bf-c93f-4a77-8128-a125cc87e042.md")
$zhcn.Hyperlinks.Add($zhcn.Range("G2"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/2926b3c98e376bb8428c6be35a5abf9cfc4f56cd/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/748403bf-c93f-4a77-8128-a125cc87e042.3377a6dd9e3bc8a1d85ec6d6e94610766108e47d.zh-cn.xlf", "", "", "748403bf-c93f-4a77-8128-a125cc87e042.3377a6dd9e3bc8a1d85ec6d6e94610766108e47d.zh-cn.xlf")
$zhcn.Hyperlinks.Add($zhcn.Range("F3"), "https://github.com/OpenLocalizationTest/oltest/blob/0b105620a0bcd802b94d27efa52426285c9e8580/e2e/89ec4140-7020-4012-9fe1-624c2b8a2ebb.md", "", "", "89ec4140-7020-4012-9fe1-624c2b8a2ebb.md")
$zhcn.Hyperlinks.Add($zhcn.Range("G3"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/2926b3c98e376bb8428c6be35a5abf9cfc4f56cd/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/89ec4140-7020-4012-9fe1-624c2b8a2ebb.494a8a5e1a4fe5743b433ae8722919d841a9c757.zh-cn.xlf", "", "", "89ec4140-7020-4012-9fe1-624c2b8a2ebb.494a8a5e1a4fe5743b433ae8722919d841a9c757.zh-cn.xlf")

# 3. de-de: update handback datetime (H2/H3) and add F/G columns
$dede.Range("H2").Value = "2016-03-21 18:47:00"
$dede.Range("H3").Value = "2016-03-21 18:47:00"

$dede.Range("F2").Value = "748403bf-c93f-4a77-8128-a125cc87e042.md"
$dede.Range("G2").Value = "748403bf-c93f-4a77-8128-a125cc87e042.3377a6dd9e3bc8a1d85ec6d6e94610766108e47d.de-de.xlf"
$dede.Range("F3").Value = "89ec4140-7020-4012-9fe1-624c2b8a2ebb.md"
$dede.Range("G3").Value = "89ec4140-7020-4012-9fe1-624c2b8a2ebb.494a8a5e1a4fe5743b433ae8722919d841a9c757.de-de.xlf"

$dede.Hyperlinks.Add($dede.Range("F2"), "https://github.com/OpenLocalizationTest/oltest/blob/0b105620a0bcd802b94d27efa52426285c9e8580/e2e/748403bf-c93f-4a77-8128-a125cc87e042.md", "", "", "748403bf-c93f-4a77-8128-a125cc87e042.md")
$dede.Hyperlinks.Add($dede.Range("G2"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/daad68d2982fb8b3705051ceedfe99655865d617/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/748403bf-c93f-4a77-8128-a125cc87e042.3377a6dd9e3bc8a1d85ec6d6e94610766108e47d.de-de.xlf", "", "", "748403bf-c93f-4a77-8128-a125cc87e042.3377a6dd9e3bc8a1d85ec6d6e94610766108e47d.de-de.xlf")
$dede.Hyperlinks.Add($dede.Range("F3"), "https://github.com/OpenLocalizationTest/oltest/blob/0b105620a0bcd802b94d27efa52426285c9e8580/e2e/89ec4140-7020-4012-9fe1-624c2b8a2ebb.md", "", "", "89ec4140-7020-4012-9fe1-624c2b8a2ebb.md")
$dede.Hyperlinks.Add($dede.Range("G3"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/daad68d2982fb8b3705051ceedfe99655865d617/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/89ec4140-7020-4012-9fe1-624c2b8a2ebb.494a8a5e1a4fe5743b433ae8722919d841a9c757.de-de.xlf", "", "", "89ec4140-7020-4012-9fe1-624c2b8a2ebb.494a8a5e1a4fe5743b433ae8722919d841a9c757.de-de.xlf")
